# Generate Report for Handback
# Update the localization-status workbook to reflect a completed handback:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Per-language sheets (zh-cn, de-de) get their Latest Target File /
#    Latest Handback File / Latest Handback DateTime columns populated
#  - New hyperlinks are added for the Latest Target File cells

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d3121ca11eb9242786270efe19529d36b4cb722a/e2e/"
$targetFile = "8bd8c3b3-dce7-4ee1-9407-457642999c12.md"

# --- Overview sheet: Status columns (E, F) for rows 2 and 3 ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$zh.Range("I2").Value = $targetFile
$zh.Hyperlinks.Add($zh.Range("I2"), $ghBase + $targetFile, "", "", $targetFile)
$zh.Range("J2").Value = "8bd8c3b3-dce7-4ee1-9407-457642999c12.5a6c0d8ad5498a2577bd30625c9578e8dedd74fb.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-30 23:08:15"

$zh.Range("I3").Value = $targetFile
$zh.Hyperlinks.Add($zh.Range("I3"), $ghBase + $targetFile, "", "", $targetFile)
$zh.Range("J3").Value = "8bd8c3b3-dce7-4ee1-9407-457642999c12.5a6c0d8ad5498a2577bd30625c9578e8dedd74fb.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-30 23:08:15"

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

$de.Range("I2").Value = $targetFile
$de.Hyperlinks.Add($de.Range("I2"), $ghBase + $targetFile, "", "", $targetFile)
$de.Range("J2").Value = "8bd8c3b3-dce7-4ee1-9407-457642999c12.5a6c0d8ad5498a2577bd30625c9578e8dedd74fb.de-de.xlf"
$de.Range("K2").Value = "2016-08-30 23:08:22"

$de.Range("I3").Value = $targetFile
$de.Hyperlinks.Add($de.Range("I3"), $ghBase + $targetFile, "", "", $targetFile)
$de.Range("J3").Value = "8bd8c3b3-dce7-4ee1-9407-457642999c12.5a6c0d8ad5498a2577bd30625c9578e8dedd74fb.de-de.xlf"
$de.Range("K3").Value = "2016-08-30 23:08:22"

# --- Column width adjustments (best-effort match of the wider "Status" /
#     long-filename columns after the new, longer text is in place) ---
$ov.Columns.Item(5).ColumnWidth = 29.166666666666668
$ov.Columns.Item(6).ColumnWidth = 29.166666666666668

$zh.Columns.Item(3).ColumnWidth = 29.166666666666668
$zh.Columns.Item(9).ColumnWidth = 39.166666666666664
$zh.Columns.Item(10).ColumnWidth = 39.166666666666664

$de.Columns.Item(3).ColumnWidth = 29.166666666666668
$de.Columns.Item(9).ColumnWidth = 39.166666666666664
$de.Columns.Item(10).ColumnWidth = 39.166666666666664
